$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCommitment")

# Fixed FCF import setup: add three new "CF" (cash-flow) columns (S:U)
# with header labels CF1/CF2/CF3, and sequential row numbers (1-8) in
# the first of the three columns (S2:S9).
$ws.Range("S1").Value = "CF1"
$ws.Range("T1").Value = "CF2"
$ws.Range("U1").Value = "CF3"

$values = @(1, 2, 3, 4, 5, 6, 7, 8)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 19).Value = $values[$i]
}

# Move the viewport so the new columns are visible and select the cell
# just below the last data row, matching the author's end-state.
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$ws.Range("S10").Select()
